$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AilenJourney")

# Insert a new row at 35 so existing row 35 (R-1) shifts down to row 36.
$ws.Rows.Item(35).Insert()

# Row 31: R-5 text update (tag/note stay the same)
$ws.Range("B31").Value = "me re caes bien, posta. la mayoria me escribe cualquier cosa pero con vos se puede hablar normal y eso esta buenisimo"

# Row 32: R-4 text update + new note
$ws.Range("B32").Value = "y aparte de bancarme a mi, que mas haces en tu tiempo libre? jaja"
$ws.Range("C32").Value = "Humor directo. Pregunta abierta."

# Row 33: tag renamed R-3 -> R-3b, new text + note
$ws.Range("A33").Value = "R-3b"
$ws.Range("B33").Value = "ahora estoy re enganchada con una serie nueva, me pasa siempre"
$ws.Range("C33").Value = "Enviar inmediatamente despues de R-3a, sin esperar respuesta."

# Row 34: tag renamed R-2 -> R-3a, new text + note
$ws.Range("A34").Value = "R-3a"
$ws.Range("B34").Value = "yo soy re nocturna, me quedo viendo anime hasta que no doy mas jaja"
$ws.Range("C34").Value = "Compartir algo personal. Si menciona un anime, conectar."

# Row 35 (new row): tag R-2, new text + note
$ws.Range("A35").Value = "R-2"
$ws.Range("B35").Value = "jaja re copado. y que onda, sos de los que se quedan despiertos hasta las 3am o de los que se levantan temprano?"
$ws.Range("C35").Value = "Mantener tono jugueton. Pregunta casual para conocerlo."

# Row 36 (shifted from old row 35): tag R-1, new text + note
$ws.Range("A36").Value = "R-1"
$ws.Range("B36").Value = "eyyy que onda 😊 recien te vi por aca y me diste curiosidad, contame algo de vos"
$ws.Range("C36").Value = "Agregar su nombre si se conoce. NUNCA decir 'baby' o 'babe'."
